$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.465.65"
$ws.Range("E2").Value = "  -3.14%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.657.13"
$ws.Range("E3").Value = "  -3.33%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.26"
$ws.Range("E5").Value = "  -1.92%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.511"
$ws.Range("E6").Value = "  -2.04%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "24.31"
$ws.Range("E8").Value = "  +2.52%  "
$ws.Range("E9").Value = "  -1.24%  "
$ws.Range("E10").Value = "  -1.99%  "
$ws.Range("E11").Value = "  -1.57%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.890.51"
$ws.Range("E12").Value = "  -3.73%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.662.93"
$ws.Range("E13").Value = "  -3.20%  "
$ws.Range("E14").Value = "  -1.81%  "
$ws.Range("E15").Value = "  +1.38%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.77"
$ws.Range("E16").Value = "  -2.33%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "27.481.12"
$ws.Range("E17").Value = "  -3.00%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "235.98"
$ws.Range("E18").Value = "  -4.34%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0₃0728"
$ws.Range("E19").Value = "  -2.32%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.67"
$ws.Range("E20").Value = "  -0.73%  "
$ws.Range("E21").Value = "  +0.10%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.44"
$ws.Range("E22").Value = "  -2.50%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.32"
$ws.Range("E23").Value = "  -2.84%  "
$ws.Range("E24").Value = "  -0.99%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.19"
$ws.Range("E25").Value = "  -1.69%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.20"
$ws.Range("E26").Value = "  -2.47%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.19"
$ws.Range("E27").Value = "  -1.73%  "
$ws.Range("E28").Value = "  -0.05%  "
$ws.Range("E29").Value = "  -1.93%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0498"
$ws.Range("E30").Value = "  -2.41%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.19"
$ws.Range("E31").Value = "  -0.76%  "
$ws.Range("E32").Value = "  -2.43%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.450.85"
$ws.Range("E33").Value = "  -1.55%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.08"
$ws.Range("E34").Value = "  -4.06%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.56"
$ws.Range("E35").Value = "  -3.66%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.40"
$ws.Range("E36").Value = "  -0.73%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.918"
$ws.Range("E37").Value = "  -5.01%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.571"
$ws.Range("E38").Value = "  -4.15%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0170"
$ws.Range("E39").Value = "  -2.84%  "
$ws.Range("E40").Value = "  -0.35%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.00"
$ws.Range("E41").Value = "  -0.09%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "66.38"
$ws.Range("E42").Value = "  -4.17%  "
$ws.Range("E43").Value = "  -3.02%  "
$ws.Range("E44").Value = "  -2.37%  "
$ws.Range("B45").Value = "TrustWalletToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.788"
$ws.Range("E45").Value = "  -1.64%  "
$ws.Range("B46").Value = "RocketPoolETH"
$ws.Range("C46").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.798.89"
$ws.Range("E46").Value = "  -3.73%  "
$ws.Range("E47").Value = "  +0.34%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "88.75"
$ws.Range("E48").Value = "  -1.49%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0₆0106"
$ws.Range("E49").Value = "  -2.13%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.101"
$ws.Range("E50").Value = "  -1.50%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.81"
$ws.Range("E51").Value = "  -3.21%  "
